$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of test data for editing a post (mirrors the existing "create post" rows)
$ws.Range("A6").Value = "EditPostSuccessfuly"
$ws.Range("B6").Value = "Test Post Edit"
$ws.Range("C6").Value = "Test first Test Edition"

$ws.Range("A7").Value = "EditPostWithoutTitle"
$ws.Range("B7").Value = " "
$ws.Range("C7").Value = "Test first Test Edition"

$ws.Range("A8").Value = "EditPostWithoutContent"
$ws.Range("B8").Value = "Test Post Edit"
$ws.Range("C8").Value = " "

$ws.Range("A9").Value = "CancelButtonEditPost"
$ws.Range("B9").Value = "Test Post"
$ws.Range("C9").Value = "Test first Test Edition"

# Update the selected cell to match the new view state from the source workbook
$ws.Range("C17").Select()
